$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting existing rows 97:225 down to 98:226
$ws.Rows("97:97").Insert()

# Populate the new row 97 with the new weekly data entry
$ws.Range("A97").Value = 10
$ws.Range("B97").Value = "Vega Modelo de Temuco"
$ws.Range("C97").Value = "La Araucanía"
$ws.Range("D97").Value = 44546
$ws.Range("E97").Value = 9
$ws.Range("F97").Value = 100112017
$ws.Range("G97").Value = "Apio"
$ws.Range("H97").Value = "Americana (o)"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 95
$ws.Range("K97").Value = 9000
$ws.Range("L97").Value = 9000
$ws.Range("M97").Value = 9000
$ws.Range("N97").Value = '$/docena de matas'
$ws.Range("O97").Value = "Provincia de Limarí"
$ws.Range("P97").Value = 1500
$ws.Range("Q97").Value = 6
$ws.Range("R97").Value = "Hortaliza"
